# "update the data analysis" - reposition the 5 shapes that make up the
# "2. Accuracy in Adam optimizer" data-analysis cluster on slide 3
# (label, plot + its two "Learning rate" callouts, and their companion
# label) - all moved together by the same offset
# (dx = +1057436 EMU, dy = -273378 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $cand = $slide.Shapes.Item($i)
        if ($cand.Id -eq $id) {
            return $cand
        }
    }
    return $null
}

# New Left/Top values below are expressed in points (EMU / 12700), chosen
# so that PowerPoint's internal Single-precision storage round-trips back
# to the exact target EMU offsets from the authoritative edit:
#   id24: (6817179,1502228) -> (7874615,1228850)
#   id25: (5070765,4221854) -> (6128201,3948476)
#   id26: (7988673,5959441) -> (9046109,5686063)
#   id27: (5909723,1898293) -> (6967159,1624915)
#   id29: (7988673,3569914) -> (9046109,3296536)

$sh = Get-ShapeById $s 24
$sh.Left = 620.0484619140625
$sh.Top = 96.75984954833984

$sh = Get-ShapeById $s 25
$sh.Left = 482.5355224609375
$sh.Top = 310.90362548828125

$sh = Get-ShapeById $s 26
$sh.Left = 712.2920532226562
$sh.Top = 447.72149658203125

$sh = Get-ShapeById $s 27
$sh.Left = 548.59521484375
$sh.Top = 127.9460678100586

$sh = Get-ShapeById $s 29
$sh.Left = 712.2920532226562
$sh.Top = 259.5697937011719
